$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

Set-TextValue 'D2' '26.441.30'
Set-TextValue 'E2' '  +1.53%  '
Set-TextValue 'D3' '1.692.08'
Set-TextValue 'E3' '  +1.52%  '
Set-TextValue 'E4' '  +0.54%  '
Set-TextValue 'D5' '219.03'
Set-TextValue 'E5' '  +1.48%  '
Set-TextValue 'D6' '0.5532'
Set-TextValue 'E6' '  +8.61%  '
Set-TextValue 'E7' '  +0.49%  '
Set-TextValue 'D8' '0.2721'
Set-TextValue 'D9' '0.06495'
Set-TextValue 'E9' '  +1.59%  '
Set-TextValue 'D10' '22.13'
Set-TextValue 'E10' '  +1.17%  '
Set-TextValue 'D11' '0.07602'
Set-TextValue 'E11' '  +1.88%  '
Set-TextValue 'B12' 'WrappedEther'
Set-TextValue 'C12' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D12' '1.695.63'
Set-TextValue 'E12' '  +1.46%  '
Set-TextValue 'B13' 'Polkadot'
Set-TextValue 'C13' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D13' '4.567'
Set-TextValue 'E13' '  +1.24%  '
Set-TextValue 'D14' '0.5854'
Set-TextValue 'D15' '0.000008472'
Set-TextValue 'E15' '  -0.30%  '
Set-TextValue 'D16' '65.49'
Set-TextValue 'E16' '  +2.29%  '
Set-TextValue 'D17' '26.525.36'
Set-TextValue 'E17' '  +1.52%  '
Set-TextValue 'D18' '4.973'
Set-TextValue 'E18' '  +1.13%  '
Set-TextValue 'E20' '  +2.08%  '
Set-TextValue 'D21' '190.93'
Set-TextValue 'E21' '  +0.62%  '
Set-TextValue 'D22' '6.267'
Set-TextValue 'E22' '  +1.36%  '
Set-TextValue 'D23' '1.010'
Set-TextValue 'E23' '  +0.49%  '
Set-TextValue 'D24' '150.07'
Set-TextValue 'E24' '  +3.50%  '
Set-TextValue 'D25' '0.1314'
Set-TextValue 'E25' '  +9.08%  '
Set-TextValue 'D26' '7.945'
Set-TextValue 'E26' '  +4.58%  '
Set-TextValue 'D27' '15.82'
Set-TextValue 'E27' '  +1.12%  '
Set-TextValue 'D28' '0.06326'
Set-TextValue 'E28' '  -4.33%  '
Set-TextValue 'D29' '1.409'
Set-TextValue 'E29' '  +5.99%  '
Set-TextValue 'D30' '1.331'
Set-TextValue 'E30' '  +1.37%  '
Set-TextValue 'D31' '3.591'
Set-TextValue 'E31' '  +1.13%  '
Set-TextValue 'D32' '3.593'
Set-TextValue 'E32' '  +2.22%  '
Set-TextValue 'D33' '1.682'
Set-TextValue 'E33' '  +1.45%  '
Set-TextValue 'D34' '1.048'
Set-TextValue 'E34' '  +3.24%  '
Set-TextValue 'D35' '0.6261'
Set-TextValue 'E35' '  +2.01%  '
Set-TextValue 'D36' '2.403'
Set-TextValue 'E36' '  +1.44%  '
Set-TextValue 'E37' '  +1.28%  '
Set-TextValue 'D38' '6.252'
Set-TextValue 'E38' '  -1.83%  '
Set-TextValue 'D39' '1.125.62'
Set-TextValue 'E39' '  +3.09%  '
Set-TextValue 'D40' '0.01647'
Set-TextValue 'E40' '  +3.44%  '
Set-TextValue 'E42' '  +0.69%  '
Set-TextValue 'D43' '100.87'
Set-TextValue 'E43' '  -0.43%  '
Set-TextValue 'D44' '1.843.54'
Set-TextValue 'E44' '  +1.66%  '
Set-TextValue 'D45' '0.00000000112'
Set-TextValue 'E45' '  -2.79%  '
Set-TextValue 'D46' '57.64'
Set-TextValue 'E46' '  +2.41%  '
Set-TextValue 'D47' '8.262'
Set-TextValue 'E47' '  +2.63%  '
Set-TextValue 'D48' '1.009'
Set-TextValue 'E48' '  +0.06%  '
Set-TextValue 'D49' '0.05282'
Set-TextValue 'E49' '  +1.02%  '
Set-TextValue 'E50' '  +1.80%  '
Set-TextValue 'E51' '  +0.26%  '
